$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: populate new "Role" column (C) for the existing 4 data rows first,
# so "Accounting" and "Developer" are appended to the shared-string table
# before the header strings.
$ws.Range("C1").Value = "Accounting"
$ws.Range("C2").Value = "Accounting"
$ws.Range("C3").Value = "Developer"
$ws.Range("C4").Value = "Developer"

# Step 2: insert a new header row at the top, shifting existing data down.
$ws.Rows.Item(1).Insert()

# Step 3: fill the new header row, right-to-left (C, B, A) so the shared
# strings get appended in the order Role, Hours, Name.
$ws.Range("C1").Value = "Role"
$ws.Range("B1").Value = "Hours"
$ws.Range("A1").Value = "Name"

# Update selection to match target state
$ws.Range("G5").Select()
